$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2810.6572
$ws.Range("I62").Value = 2326.1924
$ws.Range("K62").Value = 2326.1924
$ws.Range("M62").Value = -1702.1924
$ws.Range("H65").Value = 2810.6572
$ws.Range("I65").Value = 2326.1924
$ws.Range("K65").Value = 11630.962
$ws.Range("M65").Value = -8510.962
$ws.Range("H76").Value = 3249.9
$ws.Range("I76").Value = 3199.75
$ws.Range("J76").Value = 3283.3333
$ws.Range("K76").Value = 3199.75
$ws.Range("L76").Value = 3283.3333
$ws.Range("M76").Value = -2884.75
$ws.Range("N76").Value = -3913.3333
$ws.Range("H79").Value = 3249.9
$ws.Range("I79").Value = 3199.75
$ws.Range("J79").Value = 3283.3333
$ws.Range("K79").Value = 3199.75
$ws.Range("L79").Value = 3283.3333
$ws.Range("M79").Value = -2107.75
$ws.Range("N79").Value = -5467.3333
$ws.Range("H98").Value = 592.6429000000001
$ws.Range("I98").Value = 309.7
$ws.Range("J98").Value = 1300
$ws.Range("K98").Value = 309.7
$ws.Range("L98").Value = 1300
$ws.Range("M98").Value = 1188.3
$ws.Range("N98").Value = -4296
$ws.Range("H122").Value = 592.6429000000001
$ws.Range("I122").Value = 309.7
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 929.0999999999999
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = 1520.9
$ws.Range("N122").Value = -8800
$ws.Range("H129").Value = 265284.1
$ws.Range("I129").Value = 495
$ws.Range("J129").Value = 279994.6
$ws.Range("K129").Value = 1485
$ws.Range("L129").Value = 839983.7999999999
$ws.Range("M129").Value = 3515
$ws.Range("N129").Value = -849983.7999999999
$ws.Range("H132").Value = 2300.7104
$ws.Range("I132").Value = 2602.6365
$ws.Range("J132").Value = 308
$ws.Range("K132").Value = 7807.9095
$ws.Range("L132").Value = 924
$ws.Range("M132").Value = -5277.9095
$ws.Range("N132").Value = -5984
$ws.Range("H135").Value = 13517197
$ws.Range("I135").Value = 633.3214
$ws.Range("J135").Value = 55568730
$ws.Range("K135").Value = 5699.8926
$ws.Range("L135").Value = 500118570
$ws.Range("M135").Value = -3164.8926
$ws.Range("N135").Value = -500123640
$ws.Range("H137").Value = 41872.08
$ws.Range("I137").Value = 1246.0769
$ws.Range("J137").Value = 85883.586
$ws.Range("K137").Value = 3738.2307
$ws.Range("L137").Value = 257650.758
$ws.Range("M137").Value = -1188.2307
$ws.Range("N137").Value = -262750.758
$ws.Range("H141").Value = 1628.5428
$ws.Range("I141").Value = 1130.2174
$ws.Range("J141").Value = 2583.6667
$ws.Range("K141").Value = 3390.6522
$ws.Range("L141").Value = 7751.000100000001
$ws.Range("M141").Value = 1789.3478
$ws.Range("N141").Value = -18111.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2801.125
$ws.Range("I2").Value = 2652.75
$ws.Range("J2").Value = 2949.5
$ws.Range("K2").Value = 2652.75
$ws.Range("L2").Value = 2949.5
$ws.Range("M2").Value = -2539.75
$ws.Range("N2").Value = -3175.5
$ws.Range("H61").Value = 900.1087
$ws.Range("I61").Value = 912.2619
$ws.Range("J61").Value = 772.5
$ws.Range("K61").Value = 912.2619
$ws.Range("L61").Value = 772.5
$ws.Range("M61").Value = -700.2619
$ws.Range("N61").Value = -1196.5
$ws.Range("H74").Value = 18519832
$ws.Range("I74").Value = 21278026
$ws.Range("K74").Value = 21278026
$ws.Range("M74").Value = -21277152
$ws.Range("H77").Value = 18519832
$ws.Range("I77").Value = 21278026
$ws.Range("K77").Value = 106390130
$ws.Range("M77").Value = -106385762
$ws.Range("H116").Value = 2801.125
$ws.Range("I116").Value = 2652.75
$ws.Range("J116").Value = 2949.5
$ws.Range("K116").Value = 2652.75
$ws.Range("L116").Value = 2949.5
$ws.Range("M116").Value = -358.75
$ws.Range("N116").Value = -7537.5
$ws.Range("H132").Value = 26206.117
$ws.Range("I132").Value = 1207.7
$ws.Range("K132").Value = 3623.1
$ws.Range("M132").Value = -1093.1
$ws.Range("H136").Value = 900.1087
$ws.Range("I136").Value = 912.2619
$ws.Range("J136").Value = 772.5
$ws.Range("K136").Value = 2736.7857
$ws.Range("L136").Value = 2317.5
$ws.Range("M136").Value = -186.7856999999999
$ws.Range("N136").Value = -7417.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2801.125
$ws.Range("I3").Value = 2652.75
$ws.Range("J3").Value = 2949.5
$ws.Range("K3").Value = 2652.75
$ws.Range("L3").Value = 2949.5
$ws.Range("M3").Value = -2538.75
$ws.Range("N3").Value = -3177.5
$ws.Range("H20").Value = 1971.6364
$ws.Range("I20").Value = 2261
$ws.Range("J20").Value = 1200
$ws.Range("K20").Value = 2261
$ws.Range("L20").Value = 1200
$ws.Range("M20").Value = -2014
$ws.Range("N20").Value = -1694
$ws.Range("H86").Value = 1487.5807
$ws.Range("I86").Value = 1373.409
$ws.Range("J86").Value = 1766.6666
$ws.Range("K86").Value = 1373.409
$ws.Range("L86").Value = 1766.6666
$ws.Range("M86").Value = -250.4090000000001
$ws.Range("N86").Value = -4012.6666
$ws.Range("H89").Value = 1487.5807
$ws.Range("I89").Value = 1373.409
$ws.Range("J89").Value = 1766.6666
$ws.Range("K89").Value = 6867.045
$ws.Range("L89").Value = 8833.333000000001
$ws.Range("M89").Value = -1251.045
$ws.Range("N89").Value = -20065.333
$ws.Range("H134").Value = 26159.756
$ws.Range("I134").Value = 31188.893
$ws.Range("J134").Value = 2900
$ws.Range("K134").Value = 93566.679
$ws.Range("L134").Value = 8700
$ws.Range("M134").Value = -91031.679
$ws.Range("N134").Value = -13770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10028.707
$ws.Range("I31").Value = 12091.827
$ws.Range("J31").Value = 5042.8335
$ws.Range("K31").Value = 12091.827
$ws.Range("L31").Value = 5042.8335
$ws.Range("M31").Value = -11796.827
$ws.Range("N31").Value = -5632.8335
$ws.Range("H34").Value = 10028.707
$ws.Range("I34").Value = 12091.827
$ws.Range("J34").Value = 5042.8335
$ws.Range("K34").Value = 12091.827
$ws.Range("L34").Value = 5042.8335
$ws.Range("M34").Value = -11889.827
$ws.Range("N34").Value = -5446.8335
$ws.Range("H132").Value = 11780.2
$ws.Range("I132").Value = 13283.167
$ws.Range("J132").Value = 3889.625
$ws.Range("K132").Value = 39849.501
$ws.Range("L132").Value = 11668.875
$ws.Range("M132").Value = -37319.501
$ws.Range("N132").Value = -16728.875
$ws.Range("H134").Value = 1115.3103
$ws.Range("I134").Value = 975
$ws.Range("J134").Value = 1344.909
$ws.Range("K134").Value = 2925
$ws.Range("L134").Value = 4034.727
$ws.Range("M134").Value = -390
$ws.Range("N134").Value = -9104.727000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1971.3334
$ws.Range("I75").Value = 1999.5
$ws.Range("J75").Value = 1915
$ws.Range("K75").Value = 5998.5
$ws.Range("L75").Value = 5745
$ws.Range("M75").Value = -5000.5
$ws.Range("N75").Value = -7741
$ws.Range("H78").Value = 1971.3334
$ws.Range("I78").Value = 1999.5
$ws.Range("J78").Value = 1915
$ws.Range("K78").Value = 17995.5
$ws.Range("L78").Value = 17235
$ws.Range("M78").Value = -13003.5
$ws.Range("N78").Value = -27219
$ws.Range("H131").Value = 795.55
$ws.Range("J131").Value = 822.3684
$ws.Range("L131").Value = 2467.1052
$ws.Range("N131").Value = -12547.1052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5777778
$ws.Range("J11").Value = 8000000
$ws.Range("L11").Value = 8000000
$ws.Range("N11").Value = -8000278
$ws.Range("H70").Value = 17357.143
$ws.Range("I70").Value = 15375
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 15375
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = -15105
$ws.Range("N70").Value = -20540
$ws.Range("H73").Value = 17357.143
$ws.Range("I73").Value = 15375
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 15375
$ws.Range("L73").Value = 20000
$ws.Range("M73").Value = -14439
$ws.Range("N73").Value = -21872
$ws.Range("H102").Value = 38465190
$ws.Range("I102").Value = 38465190
$ws.Range("K102").Value = 38465190
$ws.Range("M102").Value = -38463568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 737
$ws.Range("J46").Value = 300
$ws.Range("L46").Value = 300
$ws.Range("N46").Value = -676
$ws.Range("H82").Value = 2657.1428
$ws.Range("I82").Value = 3125
$ws.Range("J82").Value = 2033.3334
$ws.Range("K82").Value = 3125
$ws.Range("L82").Value = 2033.3334
$ws.Range("M82").Value = -2764
$ws.Range("N82").Value = -2755.3334
$ws.Range("H85").Value = 2657.1428
$ws.Range("I85").Value = 3125
$ws.Range("J85").Value = 2033.3334
$ws.Range("K85").Value = 3125
$ws.Range("L85").Value = 2033.3334
$ws.Range("M85").Value = -1877
$ws.Range("N85").Value = -4529.3334
$ws.Range("H100").Value = 2114.75
$ws.Range("I100").Value = 1515
$ws.Range("K100").Value = 1515
$ws.Range("M100").Value = -974
$ws.Range("H122").Value = 1156461.2
$ws.Range("I122").Value = 3270740.5
$ws.Range("K122").Value = 9812221.5
$ws.Range("M122").Value = -9809771.5
$ws.Range("H132").Value = 1114
$ws.Range("I132").Value = 1017.28204
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3051.84612
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -521.8461200000002
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 16188.788
$ws.Range("I136").Value = 18097.068
$ws.Range("J136").Value = 2353.75
$ws.Range("K136").Value = 54291.204
$ws.Range("L136").Value = 7061.25
$ws.Range("M136").Value = -51741.204
$ws.Range("N136").Value = -12161.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1244.08
$ws.Range("I122").Value = 1220.9166
$ws.Range("K122").Value = 3662.7498
$ws.Range("M122").Value = -1212.7498
$ws.Range("H132").Value = 881.43335
$ws.Range("I132").Value = 581.7143
$ws.Range("K132").Value = 1745.1429
$ws.Range("M132").Value = 784.8571000000002
$ws.Range("H136").Value = 28572908
$ws.Range("I136").Value = 33334660
$ws.Range("J136").Value = 2399
$ws.Range("K136").Value = 100003980
$ws.Range("L136").Value = 7197
$ws.Range("M136").Value = -100001430
$ws.Range("N136").Value = -12297
